# "Added testing to project"
#
# Adds a new row for "Sravya" (the student running a test submission) to the
# roster on Sheet1, and turns the former "Email"/"Name" header row into a
# second live entry for that same student (email hyperlink + initial),
# mirroring how the sheet's own header row doubles as a data row in the
# source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: Sravya's entry ------------------------------------------
$ws.Range("A6").Value = "s@nwmissouri.edu"
$ws.Range("B6").Value = "Sravya"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:s@nwmissouri.edu")
$ws.Range("A6").Style = "Hyperlink"

# --- Row 1 (former header) becomes a second entry for Sravya -----------
# Set B1 before A1 so the shared-string table's insertion order matches.
$ws.Range("B1").Value = "S"
$ws.Range("A1").Value = "Email@nwmissouri.edu"
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:Email@nwmissouri.edu")
$ws.Range("A1").Style = "Hyperlink"

# --- Move the active selection, as left by the editing session ---------
$ws.Range("D7").Select() | Out-Null
